$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1 and F2) - set in shared-string insertion order
$ws.Range("F2").Value = "Valores son avg r from P1 perspective"
$ws.Range("B1").Value = "MODEL1"

# Column B fill rows 3-52 with iteration counts
$ws.Range("B3").Value = 1000000
$ws.Range("B4").Value = 2000000
$ws.Range("B5").Value = 3000000
$ws.Range("B6").Value = 4000000
$ws.Range("B7").Value = 5000000
$ws.Range("B8").Value = 6000000
$ws.Range("B9").Value = 7000000
$ws.Range("B10").Value = 8000000
$ws.Range("B11").Value = 9000000
$ws.Range("B12").Value = 10000000
$ws.Range("B13").Value = 11000000
$ws.Range("B14").Value = 12000000
$ws.Range("B15").Value = 13000000
$ws.Range("B16").Value = 14000000
$ws.Range("B17").Value = 15000000
$ws.Range("B18").Value = 16000000
$ws.Range("B19").Value = 17000000
$ws.Range("B20").Value = 18000000
$ws.Range("B21").Value = 19000000
$ws.Range("B22").Value = 20000000
$ws.Range("B23").Value = 21000000
$ws.Range("B24").Value = 22000000
$ws.Range("B25").Value = 23000000
$ws.Range("B26").Value = 24000000
$ws.Range("B27").Value = 25000000
$ws.Range("B28").Value = 26000000
$ws.Range("B29").Value = 27000000
$ws.Range("B30").Value = 28000000
$ws.Range("B31").Value = 29000000
$ws.Range("B32").Value = 30000000
$ws.Range("B33").Value = 31000000
$ws.Range("B34").Value = 32000000
$ws.Range("B35").Value = 33000000
$ws.Range("B36").Value = 34000000
$ws.Range("B37").Value = 35000000
$ws.Range("B38").Value = 36000000
$ws.Range("B39").Value = 37000000
$ws.Range("B40").Value = 38000000
$ws.Range("B41").Value = 39000000
$ws.Range("B42").Value = 40000000
$ws.Range("B43").Value = 41000000
$ws.Range("B44").Value = 42000000
$ws.Range("B45").Value = 43000000
$ws.Range("B46").Value = 44000000
$ws.Range("B47").Value = 45000000
$ws.Range("B48").Value = 46000000
$ws.Range("B49").Value = 47000000
$ws.Range("B50").Value = 48000000
$ws.Range("B51").Value = 49000000
$ws.Range("B52").Value = 50000000

# Columns C and D updated/new values for specific rows
$ws.Range("C3").Value = -1.79
$ws.Range("D3").Value = 3.12
$ws.Range("C4").Value = -1.15
$ws.Range("D4").Value = 1.28
$ws.Range("C5").Value = -0.57
$ws.Range("D5").Value = 0.56
$ws.Range("C6").Value = -0.46
$ws.Range("D6").Value = 0.59
$ws.Range("C7").Value = -0.44
$ws.Range("D7").Value = 0.34
$ws.Range("C8").Value = -0.33
$ws.Range("D8").Value = 0.28
$ws.Range("C9").Value = -0.37
$ws.Range("D9").Value = 0.29
$ws.Range("C10").Value = -0.38
$ws.Range("D10").Value = 0.28
$ws.Range("C11").Value = -1.38
$ws.Range("D11").Value = 0.48
$ws.Range("C12").Value = -0.5
$ws.Range("D12").Value = 0.44
$ws.Range("C22").Value = -0.48
$ws.Range("D22").Value = 0.86
$ws.Range("C46").Value = -0.29
$ws.Range("D46").Value = 0.31

# Update selection to match target view state
$ws.Range("B2").Select() | Out-Null
